$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$rows = @(
    @{Row=2; D="41.832.37"; E="  -0.29%  "},
    @{Row=3; D="2.264.29"; E="  -0.47%  "},
    @{Row=4; E="  -0.06%  "},
    @{Row=5; D="303.44"; E="  +0.34%  "},
    @{Row=6; D="92.62"; E="  -0.33%  "},
    @{Row=7; D="0.530"; E="  +0.96%  "},
    @{Row=8; E="  -0.05%  "},
    @{Row=9; D="0.485"; E="  -0.67%  "},
    @{Row=10; D="32.44"; E="  -0.71%  "},
    @{Row=11; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.0797"; E="  -0.40%  "},
    @{Row=12; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.113"; E="  -1.91%  "},
    @{Row=13; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="6.67"; E="  -0.50%  "},
    @{Row=14; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.615.31"; E="  -0.45%  "},
    @{Row=15; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="14.28"; E="  +0.60%  "},
    @{Row=16; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.267.83"; E="  -0.10%  "},
    @{Row=17; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.781"; E="  +3.36%  "},
    @{Row=18; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="41.760.37"; E="  -0.16%  "},
    @{Row=19; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="12.67"; E="  +3.54%  "},
    @{Row=20; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.0₃0908"; E="  -0.05%  "},
    @{Row=21; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.93"; E="  -0.10%  "},
    @{Row=22; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="67.60"; E="  +0.44%  "},
    @{Row=23; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="244.19"; E="  +1.24%  "},
    @{Row=24; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="2.59"; E="  +0.16%  "},
    @{Row=25; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.93"; E="  +2.68%  "},
    @{Row=26; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.00"; E="  +0.03%  "},
    @{Row=27; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="24.01"; E="  +0.70%  "},
    @{Row=28; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="9.59"; E="  -1.27%  "},
    @{Row=29; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.07"; E="  -5.62%  "},
    @{Row=30; B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="34.95"; E="  +2.52%  "},
    @{Row=31; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="160.05"; E="  +0.86%  "},
    @{Row=32; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.29"; E="  +1.92%  "},
    @{Row=33; B="FirstDigitalUSD"; C="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D="0.999"; E="  -0.13%  "},
    @{Row=34; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.0744"; E="  +0.65%  "},
    @{Row=35; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="3.01"; E="  -1.69%  "},
    @{Row=36; B="Celestia"; C="https://coinranking.com/coin/YQcD0lBl7+celestia-tia"; D="16.93"; E="  +1.99%  "},
    @{Row=37; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.106"; E="  +1.47%  "},
    @{Row=38; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="2.37"; E="  -1.11%  "},
    @{Row=39; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.116"; E="  +1.05%  "},
    @{Row=40; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.80"; E="  +0.04%  "},
    @{Row=41; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="3.92"; E="  -1.72%  "},
    @{Row=42; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="19.97"; E="  -1.26%  "},
    @{Row=43; D="2.007.01"; E="  -2.60%  "},
    @{Row=44; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0282"; E="  +1.41%  "},
    @{Row=45; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="10.33"; E="  +1.62%  "},
    @{Row=46; B="ApeXProtocol"; C="https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D="2.16"; E="  +7.93%  "},
    @{Row=47; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="2.90"; E="  -1.68%  "},
    @{Row=48; B="MultiversX"; C="https://coinranking.com/coin/omwkOTglq+multiversx-egld"; D="52.99"; E="  +2.48%  "},
    @{Row=49; B="BitcoinSV"; C="https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; D="73.23"; E="  +3.48%  "},
    @{Row=50; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.15"; E="  +0.06%  "},
    @{Row=51; D="1.50"; E="  -0.34%  "}
)

foreach ($r in $rows) {
    if ($r.ContainsKey("B")) { Set-TextValue $ws.Range("B" + $r.Row) $r.B }
    if ($r.ContainsKey("C")) { Set-TextValue $ws.Range("C" + $r.Row) $r.C }
    if ($r.ContainsKey("D")) { Set-TextValue $ws.Range("D" + $r.Row) $r.D }
    if ($r.ContainsKey("E")) { Set-TextValue $ws.Range("E" + $r.Row) $r.E }
}
